# The deck currently carries two DrawingML theme parts:
#   ppt/theme/theme1.xml -> bound to the (only) slide master -> "Integral" / "Red Violet"
#   ppt/theme/theme2.xml -> bound only to the notes master   -> "Office Theme" / "Office"
#
# The target edit swaps their contents (the stock "Office Theme" palette becomes
# the live design, the old "Integral"/Red-Violet palette is pushed onto the
# notes-only theme part), with font scheme / format scheme untouched (they were
# already identical between the two parts).
#
# PowerPoint's object model only exposes the *active* 12-slot theme color scheme
# through Slide/SlideRange.ThemeColorScheme (it is shared by every slide because
# they all hang off the single slide master/theme), so we drive the swap by
# writing the stock Office palette into that scheme, in clrScheme document order
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$officeThemeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $hex = $officeThemeColors[$i]
    $r = ($hex -shr 16) -band 0xFF
    $g = ($hex -shr 8) -band 0xFF
    $b = $hex -band 0xFF
    # VBA-style RGB() long: R + G*256 + B*65536
    $tcs.Colors($i + 1).RGB = $r + ($g * 256) + ($b * 65536)
}
